$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) are stored as text in the source data
# (prices use "." as thousands separators, e.g. "26.928.69", and percentages
# keep surrounding whitespace), so force text format before writing values to
# prevent Excel from re-interpreting numeric-looking strings as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.928.69'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '1.554.09'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = '207.09'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').Value = '0.485'
$ws.Range('E6').Value = '  +0.84%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('D8').Value = '21.71'
$ws.Range('E8').Value = '  +2.15%  '
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('D11').Value = '0.0859'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').Value = '1.775.54'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').Value = '1.561.91'
$ws.Range('E13').Value = '  +1.61%  '
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').Value = '0.515'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '61.85'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = '26.896.65'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = '215.90'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = '7.22'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = '152.38'
$ws.Range('E26').Value = '  +2.66%  '
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('E29').Value = '  +1.18%  '
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').Value = '1.417.18'
$ws.Range('E33').Value = '  +4.16%  '
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  +4.15%  '
$ws.Range('E36').Value = '  +3.40%  '
$ws.Range('E38').Value = '  +1.01%  '
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('D40').Value = '0.806'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').Value = '5.66'
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').Value = '0.988'
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').Value = '2.27'
$ws.Range('E44').Value = '  +3.93%  '
$ws.Range('D45').Value = '63.68'
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('D47').Value = '1.689.90'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('D48').Value = '86.11'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('D49').Value = '0.0518'
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').Value = '0.0960'
$ws.Range('E50').Value = '  +1.36%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.37%  '
